$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.08243182378393
$ws.Range("C2").Value = 11.54652060741781
$ws.Range("D2").Value = 14.93673672045921
$ws.Range("E2").Value = 16.34864725508803
$ws.Range("G2").Value = 3.65508773797977
$ws.Range("I2").Value = 22.16266826565719
$ws.Range("J2").Value = 9.328658768441201
$ws.Range("N2").Value = 17.20295148930182
$ws.Range("O2").Value = 26.01961115084887
$ws.Range("B3").Value = 15.56285996144461
$ws.Range("C3").Value = 11.07790983594115
$ws.Range("D3").Value = 14.8751389402851
$ws.Range("E3").Value = 16.28739504374102
$ws.Range("G3").Value = 3.657771447183607
$ws.Range("I3").Value = 22.25870033585676
$ws.Range("J3").Value = 9.337402268129123
$ws.Range("N3").Value = 17.25894361650414
$ws.Range("O3").Value = 26.04568234186899
$ws.Range("B4").Value = 15.23735797720343
$ws.Range("C4").Value = 10.78154081173902
$ws.Range("D4").Value = 14.84072283646632
$ws.Range("E4").Value = 16.25354583896944
$ws.Range("G4").Value = 3.65950606170176
$ws.Range("I4").Value = 22.32313973539977
$ws.Range("J4").Value = 9.344242886257621
$ws.Range("N4").Value = 17.29519501451462
$ws.Range("O4").Value = 26.06892182955569
$ws.Range("B5").Value = 15.10329045574808
$ws.Range("C5").Value = 10.65877703487406
$ws.Range("D5").Value = 14.82756424920418
$ws.Range("E5").Value = 16.24070688236896
$ws.Range("G5").Value = 3.660234834158326
$ws.Range("I5").Value = 22.35077116019686
$ws.Range("J5").Value = 9.347400637163284
$ws.Range("N5").Value = 17.31043961739885
$ws.Range("O5").Value = 26.08020502937777
$ws.Range("B6").Value = 15.08094909030358
$ws.Range("C6").Value = 10.63827768020671
$ws.Range("D6").Value = 14.82543186506827
$ws.Range("E6").Value = 16.23863290668958
$ws.Range("G6").Value = 3.660357171354234
$ws.Range("I6").Value = 22.35544204786365
$ws.Range("J6").Value = 9.347947332840699
$ws.Range("N6").Value = 17.31299950513793
$ws.Range("O6").Value = 26.08218791198756
$ws.Range("B7").Value = 15.23555537072718
$ws.Range("C7").Value = 10.77989298225509
$ws.Range("D7").Value = 14.84054185537335
$ws.Range("E7").Value = 16.25336881053637
$ws.Range("G7").Value = 3.659515801398763
$ws.Range("I7").Value = 22.32350683454144
$ws.Range("J7").Value = 9.34428397418541
$ws.Range("N7").Value = 17.29539869650688
$ws.Range("O7").Value = 26.06906666637455
$ws.Range("B8").Value = 15.90474822913254
$ws.Range("C8").Value = 11.38683931850287
$ws.Range("D8").Value = 14.91479739619653
$ws.Range("E8").Value = 16.32675312776511
$ws.Range("G8").Value = 3.655995109402697
$ws.Range("I8").Value = 22.19464074036505
$ws.Range("J8").Value = 9.331368006886136
$ws.Range("N8").Value = 17.22186971929848
$ws.Range("O8").Value = 26.02709650231651
$ws.Range("B9").Value = 17.15732923405663
$ws.Range("C9").Value = 12.50157216499539
$ws.Range("D9").Value = 15.08694066923666
$ws.Range("E9").Value = 16.50002585673614
$ws.Range("G9").Value = 3.649776398434019
$ws.Range("I9").Value = 21.98560044216746
$ws.Range("J9").Value = 9.317718905359957
$ws.Range("N9").Value = 17.0924818982614
$ws.Range("O9").Value = 26.0023815061072
$ws.Range("B10").Value = 18.03131195603956
$ws.Range("C10").Value = 13.2666192759822
$ws.Range("D10").Value = 15.22885942311162
$ws.Range("E10").Value = 16.64453220310917
$ws.Range("G10").Value = 3.645620546633745
$ws.Range("I10").Value = 21.85894047724751
$ws.Range("J10").Value = 9.314805989157223
$ws.Range("N10").Value = 17.00637346424427
$ws.Range("O10").Value = 26.01955672491036
$ws.Range("B11").Value = 18.41709449522034
$ws.Range("C11").Value = 13.60161851112244
$ws.Range("D11").Value = 15.29659987730468
$ws.Range("E11").Value = 16.7138358869765
$ws.Range("G11").Value = 3.643818611267965
$ws.Range("I11").Value = 21.80722982455839
$ws.Range("J11").Value = 9.315023058911695
$ws.Range("N11").Value = 16.96912977281165
$ws.Range("O11").Value = 26.03506952317938
$ws.Range("B12").Value = 18.56135776544243
$ws.Range("C12").Value = 13.72651075466447
$ws.Range("D12").Value = 15.32269259745734
$ws.Range("E12").Value = 16.74057630973548
$ws.Range("G12").Value = 3.64314892571432
$ws.Range("I12").Value = 21.78850315753387
$ws.Range("J12").Value = 9.315326600755707
$ws.Range("N12").Value = 16.95530265898437
$ws.Range("O12").Value = 26.04205136695045
$ws.Range("B13").Value = 18.53037110438662
$ws.Range("C13").Value = 13.69970167672843
$ws.Range("D13").Value = 15.31705372675932
$ws.Range("E13").Value = 16.73479546304192
$ws.Range("G13").Value = 3.643292592162537
$ws.Range("I13").Value = 21.79249816487381
$ws.Range("J13").Value = 9.315251391050866
$ws.Range("N13").Value = 16.95826830253148
$ws.Range("O13").Value = 26.0404984538324
$ws.Range("B14").Value = 18.42900036235419
$ws.Range("C14").Value = 13.61193332517692
$ws.Range("D14").Value = 15.2987378027326
$ws.Range("E14").Value = 16.71602598120914
$ws.Range("G14").Value = 3.643763262341956
$ws.Range("I14").Value = 21.80567200271611
$ws.Range("J14").Value = 9.315043598203307
$ws.Range("N14").Value = 16.96798667617708
$ws.Range("O14").Value = 26.03562173044449
$ws.Range("B15").Value = 18.36666674173589
$ws.Range("C15").Value = 13.55791431755136
$ws.Range("D15").Value = 15.28757568167724
$ws.Range("E15").Value = 16.70459330705763
$ws.Range("G15").Value = 3.64405320913496
$ws.Range("I15").Value = 21.81385287286915
$ws.Range("J15").Value = 9.31494513026208
$ws.Range("N15").Value = 16.9739754119961
$ws.Range("O15").Value = 26.03277881272181
$ws.Range("B16").Value = 18.00585151959094
$ws.Range("C16").Value = 13.24445604672194
$ws.Range("D16").Value = 15.22449503873845
$ws.Range("E16").Value = 16.64007344833746
$ws.Range("G16").Value = 3.645740083877439
$ws.Range("I16").Value = 21.86243923969766
$ws.Range("J16").Value = 9.31482280185393
$ws.Range("N16").Value = 17.00884613476908
$ws.Range("O16").Value = 26.01869791385034
$ws.Range("B17").Value = 17.7813860502783
$ws.Range("C17").Value = 13.0487540099274
$ws.Range("D17").Value = 15.18660031557105
$ws.Range("E17").Value = 16.6013948901076
$ws.Range("G17").Value = 3.646797565189925
$ws.Range("I17").Value = 21.8937624282218
$ws.Range("J17").Value = 9.315142475972904
$ws.Range("N17").Value = 17.03073121277227
$ws.Range("O17").Value = 26.01203236778381
$ws.Range("B18").Value = 17.65117578505883
$ws.Range("C18").Value = 12.93497001639143
$ws.Range("D18").Value = 15.16510471872181
$ws.Range("E18").Value = 16.57948478797758
$ws.Range("G18").Value = 3.647414142412789
$ws.Range("I18").Value = 21.91233442083524
$ws.Range("J18").Value = 9.315471531633332
$ws.Range("N18").Value = 17.04350040534816
$ws.Range("O18").Value = 26.00892330849787
$ws.Range("B19").Value = 17.60690340657994
$ws.Range("C19").Value = 12.89623796114702
$ws.Range("D19").Value = 15.15787877091424
$ws.Range("E19").Value = 16.57212472847307
$ws.Range("G19").Value = 3.647624339738014
$ws.Range("I19").Value = 21.91871785527238
$ws.Range("J19").Value = 9.315607891544101
$ws.Range("N19").Value = 17.04785503098348
$ws.Range("O19").Value = 26.00799509064076
$ws.Range("B20").Value = 17.80539598000489
$ws.Range("C20").Value = 13.06971393929403
$ws.Range("D20").Value = 15.19060329569241
$ws.Range("E20").Value = 16.60547754832951
$ws.Range("G20").Value = 3.646684131645741
$ws.Range("I20").Value = 21.89037046582897
$ws.Range("J20").Value = 9.315093422482345
$ws.Range("N20").Value = 17.02838273353168
$ws.Range("O20").Value = 26.01266690824816
$ws.Range("B21").Value = 18.45882586267414
$ws.Range("C21").Value = 13.63776700455522
$ws.Range("D21").Value = 15.30410580395732
$ws.Range("E21").Value = 16.7215256892471
$ws.Range("G21").Value = 3.643624671842124
$ws.Range("I21").Value = 21.80177927796361
$ws.Range("J21").Value = 9.315098628661794
$ws.Range("N21").Value = 16.96512466360482
$ws.Range("O21").Value = 26.03702408887069
$ws.Range("B22").Value = 18.87519783773749
$ws.Range("C22").Value = 13.99753725206908
$ws.Range("D22").Value = 15.38084782583802
$ws.Range("E22").Value = 16.80025637175305
$ws.Range("G22").Value = 3.641698946259513
$ws.Range("I22").Value = 21.74886614410391
$ws.Range("J22").Value = 9.316391919796445
$ws.Range("N22").Value = 16.92539173509194
$ws.Range("O22").Value = 26.05939751380416
$ws.Range("B23").Value = 18.65398827423311
$ws.Range("C23").Value = 13.8065989626392
$ws.Range("D23").Value = 15.33966044488275
$ws.Range("E23").Value = 16.75797782371305
$ws.Range("G23").Value = 3.642720011727442
$ws.Range("I23").Value = 21.77664881572104
$ws.Range("J23").Value = 9.315583804887801
$ws.Range("N23").Value = 16.94645093448591
$ws.Range("O23").Value = 26.04686599799426
$ws.Range("B24").Value = 17.79454470047507
$ws.Range("C24").Value = 13.06024191104762
$ws.Range("D24").Value = 15.18879264200356
$ws.Range("E24").Value = 16.60363075965066
$ws.Range("G24").Value = 3.6467353881025
$ws.Range("I24").Value = 21.8919022159133
$ws.Range("J24").Value = 9.315115147027047
$ws.Range("N24").Value = 17.02944389786452
$ws.Range("O24").Value = 26.01237778017179
$ws.Range("B25").Value = 16.82595023421456
$ws.Range("C25").Value = 12.20897755885062
$ws.Range("D25").Value = 15.03760228325613
$ws.Range("E25").Value = 16.45007496069634
$ws.Range("G25").Value = 3.651385847882736
$ws.Range("I25").Value = 22.03744537808064
$ws.Range("J25").Value = 9.320161162771551
$ws.Range("N25").Value = 17.12590749235675
$ws.Range("O25").Value = 26.00287494562365
